# PL_specific_ion_cfg.xlsx edit:
# Insert a new "CHARGE_MODE" column (POS/NEG) right before the existing
# "PR_CHARGE" column (old column E -> new column F), populate it for every
# data row based on the ion mode already present in PR_CHARGE, and move the
# active selection to the new header cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E; everything from old E onward (PR_CHARGE, LABEL,
# REMARKS) shifts one column to the right (E->F, F->G, G->H).
$ws.Columns("E:E").Insert()

# Match the column width used for the new CHARGE_MODE column.
$ws.Columns("E:E").ColumnWidth = 15.28515625

# Data rows (row 1 is the header, blank rows 4/15/22/27/28/33 are spacers).
$dataRows = 2,3,5,6,7,8,9,10,11,12,13,14,16,17,18,19,20,21,23,24,25,26,29,30,31,32,34,35,36,37,38,39

foreach ($r in $dataRows) {
    $prCharge = $ws.Range("F$r").Value()
    if ($prCharge -eq "[M+H]+" -or $prCharge -eq "[M+NH4]+") {
        $ws.Range("E$r").Value() = "POS"
    } elseif ($prCharge -eq "[M-H]-" -or $prCharge -eq "[M+HCOO]-") {
        $ws.Range("E$r").Value() = "NEG"
    }
}

# Header for the new column, set last so it becomes the newest shared string.
$ws.Range("E1").Value() = "CHARGE_MODE"

$ws.Range("E1").Select()
